# Update the maze grid values on Sheet1 (spiral_hole3) and adjust the
# selected cell to match the committed state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New grid values (10 rows x 10 cols, columns A..J)
$values = @(
    @(1,0,1,1,1,1,1,1,1,1),
    @(1,0,1,0,0,1,0,0,0,1),
    @(1,0,1,1,1,1,1,1,0,1),
    @(1,0,1,0,1,0,0,1,0,1),
    @(1,0,1,0,1,1,1,1,0,1),
    @(1,0,1,0,0,1,0,1,0,1),
    @(1,0,1,1,0,1,0,1,0,1),
    @(1,0,1,1,1,1,0,1,0,1),
    @(1,0,0,0,0,0,0,1,0,1),
    @(1,1,1,1,1,1,1,1,0,1)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Reflect the updated selection recorded in the workbook.
$ws.Range("W2").Select()

$wb.Save()
